$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New meeting block (rows 14-16), mirroring the existing "Meeting" pattern.
$ws.Range("A14").Value = "Meeting5"
$ws.Range("F14").Value = "All Member"

$ws.Range("A15").Value = "Planing for next module"
$ws.Range("F15").Value = "All Member"

$ws.Range("A16").Value = "DOM manipulating"
$ws.Range("G16").Value = "夏义"

# Match the bold style (s="6") used by the rest of column A / task rows.
$ws.Range("A14:A16").Font.Bold = $true
$ws.Range("F14:F15").Font.Bold = $true
$ws.Range("G16").Font.Bold = $true

# New column F now holds data, so it gets an explicit (auto-fit-like) width,
# mirroring the other date columns that already carry a customWidth.
$ws.Range("F1").ColumnWidth = 9.666666666666666

# Keep the active cell/selection consistent with the appended rows.
$ws.Range("A17").Select()
